# Update the "Förändrad" (Changed) date column (C) from 2025-02-02 (45690)
# to 2025-02-03 (45691) for all data rows (rows 2 through 36).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 36; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45690) {
        $cell.Value2 = 45691
    }
}
